$d = $word.ActiveDocument

# Update the date heading paragraph (first paragraph, before the table).
$d.Paragraphs.Item(1).Range.Text = "2024-04-10 Wednesday"

# Update the practice-problem table. Only rows 1, 5, 9, 13, 17 (1-based)
# contain text; the rest are blank working rows. Each row has 5 columns.
$table = $d.Tables.Item(1)

$newValues = @{
    1  = @("50÷4=12, 2", "57÷8=7, 1", "78÷6=13, 0", "19÷3=6, 1", "34÷6=5, 4")
    5  = @("15÷5=3, 0", "79÷4=19, 3", "64÷6=10, 4", "12÷4=3, 0", "14÷3=4, 2")
    9  = @("21÷6=3, 3", "95÷5=19, 0", "24÷4=6, 0", "37÷5=7, 2", "51÷8=6, 3")
    13 = @("84÷6=14, 0", "38÷3=12, 2", "70÷4=17, 2", "12÷3=4, 0", "34÷6=5, 4")
    17 = @("56÷3=18, 2", "49÷8=6, 1", "35÷4=8, 3", "36÷2=18, 0", "97÷3=32, 1")
}

foreach ($row in $newValues.Keys) {
    $values = $newValues[$row]
    for ($col = 1; $col -le $values.Length; $col++) {
        $table.Cell($row, $col).Range.Text = $values[$col - 1]
    }
}

Write-Host "Done applying replacements"
